$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
# for rows 2-16 (data rows), reflecting a reshuffle of the weekly records.

$rows = @{
  2  = @{ D = 44176; J = 10; K = 4000; L = 4000; M = 4000; P = 4000 }
  3  = @{ D = 44497; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
  4  = @{ D = 44315; J = 40; K = 4000; L = 4000; M = 4000; P = 4000 }
  5  = @{ D = 44365; J = 55; K = 5000; L = 5000; M = 5000; P = 5000 }
  6  = @{ D = 44504; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 }
  7  = @{ D = 44291; J = 35; K = 4000; L = 4000; M = 4000; P = 4000 }
  8  = @{ D = 44259; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 }
  9  = @{ D = 44508; J = 30; K = 4000; L = 4000; M = 4000; P = 4000 }
  10 = @{ D = 44498; J = 40; K = 4000; L = 4000; M = 4000; P = 4000 }
  11 = @{ D = 44316; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
  12 = @{ D = 44313; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
  13 = @{ D = 44280; J = 55; K = 4000; L = 4000; M = 4000; P = 4000 }
  14 = @{ D = 44509; J = 20; K = 4000; L = 4000; M = 4000; P = 4000 }
  15 = @{ D = 44301; J = 40; K = 3000; L = 3000; M = 3000; P = 3000 }
  16 = @{ D = 44312; J = 50; K = 4000; L = 4000; M = 4000; P = 4000 }
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Range("D$r").Value = $vals.D
  $ws.Range("J$r").Value = $vals.J
  $ws.Range("K$r").Value = $vals.K
  $ws.Range("L$r").Value = $vals.L
  $ws.Range("M$r").Value = $vals.M
  $ws.Range("P$r").Value = $vals.P
}
